$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 78, shifting the existing
# rows 78-89 down to 80-91 (weekly price-sheet update: two new records
# pushed in at the top of this variety block).
$ws.Rows.Item(78).Insert()
$ws.Rows.Item(78).Insert()

# --- New row 78: Melón / Tuna / Primera ---
$ws.Cells.Item(78, 1).Value = 1
$ws.Cells.Item(78, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(78, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(78, 4).Value = 45275
$ws.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(78, 5).Value = 15
$ws.Cells.Item(78, 6).Value = 100112027
$ws.Cells.Item(78, 7).Value = "Melón"
$ws.Cells.Item(78, 8).Value = "Tuna"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 140
$ws.Cells.Item(78, 11).Value = 20000
$ws.Cells.Item(78, 12).Value = 22000
$ws.Cells.Item(78, 13).Value = 21000
$ws.Cells.Item(78, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(78, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(78, 16).Value = 1167
$ws.Cells.Item(78, 17).Value = 18
$ws.Cells.Item(78, 18).Value = "Hortaliza"

# --- New row 79: Melón / Tuna / Segunda ---
$ws.Cells.Item(79, 1).Value = 1
$ws.Cells.Item(79, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(79, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(79, 4).Value = 45275
$ws.Cells.Item(79, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(79, 5).Value = 15
$ws.Cells.Item(79, 6).Value = 100112027
$ws.Cells.Item(79, 7).Value = "Melón"
$ws.Cells.Item(79, 8).Value = "Tuna"
$ws.Cells.Item(79, 9).Value = "Segunda"
$ws.Cells.Item(79, 10).Value = 160
$ws.Cells.Item(79, 11).Value = 24000
$ws.Cells.Item(79, 12).Value = 25000
$ws.Cells.Item(79, 13).Value = 24500
$ws.Cells.Item(79, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(79, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 16).Value = 1021
$ws.Cells.Item(79, 17).Value = 24
$ws.Cells.Item(79, 18).Value = "Hortaliza"
